# Add data for 2022-12-04
#   - Rename the sheet and update the "November (through 11-25)" label to
#     reflect the new cutoff date of 11-26.
#   - Update the November row (row 12) and Total row (row 13) with the
#     refreshed year-over-year carjacking counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet/tab to reflect the new "through" date.
$ws.Name = "Through 2022-11-26"

# Update the row label for November.
$ws.Range("A12").Value = "November (through 11-26)"

# November row (row 12): 2015 .. 2022 in columns B .. I
$ws.Range("B12").Value = 28
$ws.Range("C12").Value = 66
$ws.Range("D12").Value = 98
$ws.Range("E12").Value = 58
$ws.Range("F12").Value = 46
$ws.Range("G12").Value = 182
$ws.Range("H12").Value = 176
$ws.Range("I12").Value = 97

# Total row (row 13): 2015 .. 2022 in columns B .. I
$ws.Range("B13").Value = 286
$ws.Range("C13").Value = 552
$ws.Range("D13").Value = 808
$ws.Range("E13").Value = 673
$ws.Range("F13").Value = 528
$ws.Range("G13").Value = 1239
$ws.Range("H13").Value = 1617
$ws.Range("I13").Value = 1495
